# Journal de bord - add new session entry:
# "TP A3 du 27/02/2024 (yasmf hello_world)"
#
# Row 17 (27/02/2024, group A2-4/MPAL) already exists in the workbook.
# Row 18 is a pre-formatted, still-empty placeholder row for the same date
# that needs to be filled in with the group-A3 (FSIL) entry for that session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in row 18 with the new journal entry ---
$ws.Range("A18").Value = 45349
$ws.Range("B18").Value = "FSIL"
$ws.Range("C18").Value = "P"
$ws.Range("F18").Value = "x"
$ws.Range("G18").Value = "TP sur exemples résolus yasmf => fin travail et questions sur Hello_world.RAF : travail et questions sur all_users. "
$ws.Range("I18").Value = "Sur autoload, 8 réponses, et une restitution des résultats sans passer la phase 2."

# --- Row height adjustments (content re-wrap after the edit) ---
$ws.Rows.Item(6).RowHeight = 33
$ws.Rows.Item(15).RowHeight = 17.25
$ws.Rows.Item(16).RowHeight = 17.25
$ws.Rows.Item(17).RowHeight = 50.25
$ws.Rows.Item(18).RowHeight = 47.25
